$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new September log entry was captured, so insert a fresh row above the
# existing row 31 ("Insert" in Excel shifts rows 31-77 down to 32-78 and
# widens the used range to A1:Y78).
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row with the latest September entry.
$ws.Range("R31").Value = "share anyone axis"
$ws.Range("S31").Value = "2024-09-05 16:38:59"
